$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "driver deadtime after ls on" (tsfet_dt_off) value for raa489300 changes from 15 to 30
$ws.Range("E9").Value = 30

# Match the saved cursor/selection position left by the edit
$ws.Range("E10").Select()
